# Add POR (Proof of Residence) documents for minor residents (age 14, lang fra)
# Inserts 5 new rows right after the existing age=14 block (old row 300),
# shifting all subsequent rows down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 301 (pushes old row 301.. down to 306..)
$ws.Range("A301:A305").EntireRow.Insert()

# Populate the newly inserted rows with the new POR document rules for apptyp 14 (fra)
$newRows = @(
    @("fra", 14, "POR", "COB"),
    @("fra", 14, "POR", "CRN"),
    @("fra", 14, "POR", "DOC002"),
    @("fra", 14, "POR", "DOC007"),
    @("fra", 14, "POR", "DOC009")
)

$r = 301
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    # Copy an existing "TRUE" text cell so the value is stored as text (matching
    # the rest of the is_active column) rather than being auto-coerced to a boolean.
    $ws.Range("E2").Copy($ws.Range("E$r"))
    $r = $r + 1
}
